# Intro - added reference
#
# 1) Merge the two runs of the "forested" paragraph into a single run
#    (Find & Replace the full text with itself forces the engine to
#    normalise/merge the backing runs, matching how Word collapses
#    adjacent runs with identical formatting after an edit).
$d = $word.ActiveDocument

$oldText = "At the start of the century 41.9% of Cambodia" + [char]0x2019 + "s land area was forested, and by 2012 the total forested area had been reduced by 19.8%, equating to over 1.3 million hectares (Davis et al 2015). Only 25 other countries lost more forest than Cambodia between 2000 " + [char]0x2013 + " 2012 (Hansen et al 2013). "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2) | Out-Null

# 2) Append the new reference paragraph, preceded by three blank
#    paragraphs and followed by one blank paragraph, using raw OOXML so
#    that the proofing (spell-check) marks around "Lizcano" and the
#    xml:space="preserve" run are reproduced exactly.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParasXml = (
    "<w:p $ns/>" +
    "<w:p $ns/>" +
    "<w:p $ns/>" +
    "<w:p $ns>" +
        "<w:r><w:t xml:space=`"preserve`">McSweeney, C., New, M. &amp; </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>Lizcano</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t>, G. 2010. UNDP Climate Change Country Profiles: Cambodia. Available: http://country-profiles.geog.ox.ac.uk/ [Accessed 23/06/2020].</w:t></w:r>" +
    "</w:p>" +
    "<w:p $ns/>"
)

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($newParasXml)
